$wb = $excel.ActiveWorkbook

function Set-CellValue($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

function Clear-CellValue($ws, $addr) {
    $ws.Range($addr).ClearContents()
}

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws "H46" 2153.3333
Set-CellValue $ws "I46" 1650
Set-CellValue $ws "J46" 3160
Set-CellValue $ws "K46" 4950
Set-CellValue $ws "L46" 9480
Set-CellValue $ws "M46" -4831
Set-CellValue $ws "N46" -9718
Set-CellValue $ws "H60" 2153.3333
Set-CellValue $ws "I60" 1650
Set-CellValue $ws "J60" 3160
Set-CellValue $ws "K60" 4950
Set-CellValue $ws "L60" 9480
Set-CellValue $ws "M60" -4466
Set-CellValue $ws "N60" -10448
Set-CellValue $ws "H62" 3161.2
Set-CellValue $ws "I62" 2966.6667
Set-CellValue $ws "J62" 3453
Set-CellValue $ws "K62" 2966.6667
Set-CellValue $ws "L62" 3453
Set-CellValue $ws "M62" -2342.6667
Set-CellValue $ws "N62" -4701
Set-CellValue $ws "H65" 3161.2
Set-CellValue $ws "I65" 2966.6667
Set-CellValue $ws "J65" 3453
Set-CellValue $ws "K65" 14833.3335
Set-CellValue $ws "L65" 17265
Set-CellValue $ws "M65" -11713.3335
Set-CellValue $ws "N65" -23505
Set-CellValue $ws "H132" 17545306
Set-CellValue $ws "I132" 20834010
Set-CellValue $ws "J132" 5555.778
Set-CellValue $ws "K132" 62502030
Set-CellValue $ws "L132" 16667.334
Set-CellValue $ws "M132" -62499500
Set-CellValue $ws "N132" -21727.334
Set-CellValue $ws "H133" 59750
Set-CellValue $ws "J133" 59750
Set-CellValue $ws "L133" 59750
Set-CellValue $ws "N133" -69870
Set-CellValue $ws "H137" 728.45905
Set-CellValue $ws "I137" 635.55554
Set-CellValue $ws "J137" 862.24
Set-CellValue $ws "K137" 1906.66662
Set-CellValue $ws "L137" 2586.72
Set-CellValue $ws "M137" 643.33338
Set-CellValue $ws "N137" -7686.72

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
Set-CellValue $ws "H5" 100
Set-CellValue $ws "I5" 0
Set-CellValue $ws "J5" 100
Set-CellValue $ws "K5" 0
Set-CellValue $ws "L5" 100
Clear-CellValue $ws "M5"
Set-CellValue $ws "N5" -324
Set-CellValue $ws "H6" 3530000
Set-CellValue $ws "J6" 60000
Set-CellValue $ws "L6" 60000
Set-CellValue $ws "N6" -60346
Set-CellValue $ws "H55" 35200.2
Set-CellValue $ws "J55" 35200.2
Set-CellValue $ws "L55" 35200.2
Set-CellValue $ws "N55" -35830.2
Set-CellValue $ws "H61" 1658.5454
Set-CellValue $ws "I61" 1028.4
Set-CellValue $ws "J61" 2183.6667
Set-CellValue $ws "K61" 1028.4
Set-CellValue $ws "L61" 2183.6667
Set-CellValue $ws "M61" -816.4000000000001
Set-CellValue $ws "N61" -2607.6667
Set-CellValue $ws "H74" 1417.5385
Set-CellValue $ws "I74" 1377.6
Set-CellValue $ws "K74" 1377.6
Set-CellValue $ws "M74" -503.5999999999999
Set-CellValue $ws "H77" 1417.5385
Set-CellValue $ws "I77" 1377.6
Set-CellValue $ws "K77" 6888
Set-CellValue $ws "M77" -2520
Set-CellValue $ws "H80" 12029
Set-CellValue $ws "I80" 2100
Set-CellValue $ws "J80" 16000.6
Set-CellValue $ws "K80" 2100
Set-CellValue $ws "L80" 16000.6
Set-CellValue $ws "M80" -1102
Set-CellValue $ws "N80" -17996.6
Set-CellValue $ws "H81" 0
Set-CellValue $ws "I81" 0
Set-CellValue $ws "J81" 0
Set-CellValue $ws "K81" 0
Set-CellValue $ws "L81" 0
Clear-CellValue $ws "M81"
Clear-CellValue $ws "N81"
Set-CellValue $ws "H83" 12029
Set-CellValue $ws "I83" 2100
Set-CellValue $ws "J83" 16000.6
Set-CellValue $ws "K83" 6300
Set-CellValue $ws "L83" 48001.8
Set-CellValue $ws "M83" -1308
Set-CellValue $ws "N83" -57985.8
Set-CellValue $ws "H84" 0
Set-CellValue $ws "I84" 0
Set-CellValue $ws "J84" 0
Set-CellValue $ws "K84" 0
Set-CellValue $ws "L84" 0
Clear-CellValue $ws "M84"
Clear-CellValue $ws "N84"
Set-CellValue $ws "H87" 30000
Set-CellValue $ws "J87" 30000
Set-CellValue $ws "L87" 30000
Set-CellValue $ws "N87" -32496
Set-CellValue $ws "H88" 5712.5557
Set-CellValue $ws "I88" 4706
Set-CellValue $ws "J88" 5838.375
Set-CellValue $ws "K88" 4706
Set-CellValue $ws "L88" 5838.375
Set-CellValue $ws "M88" -4300
Set-CellValue $ws "N88" -6650.375
Set-CellValue $ws "H90" 30000
Set-CellValue $ws "J90" 30000
Set-CellValue $ws "L90" 90000
Set-CellValue $ws "N90" -102480
Set-CellValue $ws "H91" 5712.5557
Set-CellValue $ws "I91" 4706
Set-CellValue $ws "J91" 5838.375
Set-CellValue $ws "K91" 4706
Set-CellValue $ws "L91" 5838.375
Set-CellValue $ws "M91" -3302
Set-CellValue $ws "N91" -8646.375
Set-CellValue $ws "H136" 1658.5454
Set-CellValue $ws "I136" 1028.4
Set-CellValue $ws "J136" 2183.6667
Set-CellValue $ws "K136" 3085.2
Set-CellValue $ws "L136" 6551.000100000001
Set-CellValue $ws "M136" -535.2000000000003
Set-CellValue $ws "N136" -11651.0001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
Set-CellValue $ws "H4" 100
Set-CellValue $ws "I4" 0
Set-CellValue $ws "J4" 100
Set-CellValue $ws "K4" 0
Set-CellValue $ws "L4" 100
Clear-CellValue $ws "M4"
Set-CellValue $ws "N4" -330
Set-CellValue $ws "H134" 7473.7393
Set-CellValue $ws "I134" 977.55554
Set-CellValue $ws "J134" 30860
Set-CellValue $ws "K134" 2932.66662
Set-CellValue $ws "L134" 92580
Set-CellValue $ws "M134" -397.66662
Set-CellValue $ws "N134" -97650

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
Set-CellValue $ws "H31" 7044031
Set-CellValue $ws "I31" 10418179
Set-CellValue $ws "J31" 2331.1738
Set-CellValue $ws "K31" 10418179
Set-CellValue $ws "L31" 2331.1738
Set-CellValue $ws "M31" -10417884
Set-CellValue $ws "N31" -2921.1738
Set-CellValue $ws "H34" 7044031
Set-CellValue $ws "I34" 10418179
Set-CellValue $ws "J34" 2331.1738
Set-CellValue $ws "K34" 10418179
Set-CellValue $ws "L34" 2331.1738
Set-CellValue $ws "M34" -10417977
Set-CellValue $ws "N34" -2735.1738
Set-CellValue $ws "H58" 1098.5312
Set-CellValue $ws "I58" 862.86365
Set-CellValue $ws "J58" 1617
Set-CellValue $ws "K58" 862.86365
Set-CellValue $ws "L58" 1617
Set-CellValue $ws "M58" -659.86365
Set-CellValue $ws "N58" -2023
Set-CellValue $ws "H105" 1440.5333
Set-CellValue $ws "I105" 1209.8182
Set-CellValue $ws "J105" 2075
Set-CellValue $ws "K105" 1209.8182
Set-CellValue $ws "L105" 2075
Set-CellValue $ws "M105" 537.1818000000001
Set-CellValue $ws "N105" -5569
Set-CellValue $ws "H136" 1098.5312
Set-CellValue $ws "I136" 862.86365
Set-CellValue $ws "J136" 1617
Set-CellValue $ws "K136" 2588.59095
Set-CellValue $ws "L136" 4851
Set-CellValue $ws "M136" -38.59094999999979
Set-CellValue $ws "N136" -9951

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
Set-CellValue $ws "H107" 500.33334
Set-CellValue $ws "I107" 383.4
Set-CellValue $ws "J107" 583.8570999999999
Set-CellValue $ws "K107" 1150.2
Set-CellValue $ws "L107" 1751.5713
Set-CellValue $ws "M107" 769.8000000000002
Set-CellValue $ws "N107" -5591.5713
Set-CellValue $ws "H134" 5682.375
Set-CellValue $ws "I134" 5793.85
Set-CellValue $ws "J134" 5125
Set-CellValue $ws "K134" 17381.55
Set-CellValue $ws "L134" 15375
Set-CellValue $ws "M134" -12311.55
Set-CellValue $ws "N134" -25515

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
Set-CellValue $ws "H62" 0
Set-CellValue $ws "J62" 0
Set-CellValue $ws "L62" 0
Clear-CellValue $ws "N62"
Set-CellValue $ws "H65" 0
Set-CellValue $ws "J65" 0
Set-CellValue $ws "L65" 0
Clear-CellValue $ws "N65"
Set-CellValue $ws "H68" 30000
Set-CellValue $ws "J68" 30000
Set-CellValue $ws "L68" 30000
Set-CellValue $ws "N68" -31622
Set-CellValue $ws "H71" 30000
Set-CellValue $ws "J71" 30000
Set-CellValue $ws "L71" 90000
Set-CellValue $ws "N71" -98112

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
Set-CellValue $ws "H47" 5399
Set-CellValue $ws "J47" 5399
Set-CellValue $ws "L47" 5399
Set-CellValue $ws "N47" -6379
Set-CellValue $ws "H52" 5399
Set-CellValue $ws "J52" 5399
Set-CellValue $ws "L52" 5399
Set-CellValue $ws "N52" -5865
Set-CellValue $ws "H132" 2484.6333
Set-CellValue $ws "I132" 2945.8108
Set-CellValue $ws "J132" 1742.7391
Set-CellValue $ws "K132" 8837.432400000002
Set-CellValue $ws "L132" 5228.2173
Set-CellValue $ws "M132" -6307.432400000002
Set-CellValue $ws "N132" -10288.2173
Set-CellValue $ws "H136" 2575.65
Set-CellValue $ws "I136" 1955.5555
Set-CellValue $ws "K136" 5866.666499999999
Set-CellValue $ws "M136" -3316.666499999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
Set-CellValue $ws "H18" 0
Set-CellValue $ws "I18" 0
Set-CellValue $ws "J18" 0
Set-CellValue $ws "K18" 0
Set-CellValue $ws "L18" 0
Clear-CellValue $ws "M18"
Clear-CellValue $ws "N18"
Set-CellValue $ws "H75" 18500
Set-CellValue $ws "I75" 0
Set-CellValue $ws "J75" 18500
Set-CellValue $ws "K75" 0
Set-CellValue $ws "L75" 18500
Clear-CellValue $ws "M75"
Set-CellValue $ws "N75" -20372
Set-CellValue $ws "H78" 18500
Set-CellValue $ws "I78" 0
Set-CellValue $ws "J78" 18500
Clear-CellValue $ws "M78"
Set-CellValue $ws "N78" -64860
Set-CellValue $ws "H132" 25513998
Set-CellValue $ws "I132" 39064080
Set-CellValue $ws "J132" 7960.1763
Set-CellValue $ws "K132" 117192240
Set-CellValue $ws "L132" 23880.5289
Set-CellValue $ws "M132" -117189710
Set-CellValue $ws "N132" -28940.5289
Set-CellValue $ws "H135" 59857.5
Set-CellValue $ws "J135" 59857.5
Set-CellValue $ws "L135" 59857.5
Set-CellValue $ws "N135" -69997.5
